$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quality_comparison")
$c = $ws.Range("C1")
$i = $c.Interior
Write-Output "Interior type: $($i.GetType())"
$f = $c.Font
Write-Output "Font type: $($f.GetType())"
